$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.121.77"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.813.44"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.12"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("E9").Value = "  +6.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0998"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "2.077.74"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "1.817.75"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.659"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.66"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "35.075.59"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "0.0₃0791"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  +4.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +19.28%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("B32").Value = "EURNeutrino"
$ws.Range("C32").Value = "https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn"
$ws.Range("D32").Value = "3.331.73"
$ws.Range("E32").Value = "  +37.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0553"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  -8.11%  "
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "92.77"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.681"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.28"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.308.34"
$ws.Range("E41").Value = "  -2.77%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.61"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.63%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "1.991.67"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("E51").Value = "  +4.70%  "
